# Generate Report for Handback
#
# The "1f1910f2-2ea9-49f2-a4ad-68d0cddcb98c" handback row (row 7) on both the
# zh-cn and de-de sheets receives a freshly generated report: a new
# "Latest Target File" hyperlink/value, a "Latest Handback File" value, an
# updated "Latest Handback DateTime", and an "Error Detail" message warning
# that the handback commit isn't the latest one. The "Error Detail" column
# is also widened so the message is readable.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/67deddc2823b61d237b596fc789cd3d7f95eccd8/e2e/1f1910f2-2ea9-49f2-a4ad-68d0cddcb98c.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ff59497e4c3192bee9f465e9252a3128e619b66a/e2e/1f1910f2-2ea9-49f2-a4ad-68d0cddcb98c.md."
$latestUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ff59497e4c3192bee9f465e9252a3128e619b66a/e2e/1f1910f2-2ea9-49f2-a4ad-68d0cddcb98c.md"
$displayName = "1f1910f2-2ea9-49f2-a4ad-68d0cddcb98c.md"

# ---- zh-cn sheet ----
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("I7").Value = $displayName
$ws.Range("I7").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("I7"), $latestUrl, "", "", $displayName)

$ws.Range("J7").Value = "1f1910f2-2ea9-49f2-a4ad-68d0cddcb98c.3d92c1c223942b2c7e80af8189297e04f76f87e6.zh-cn.xlf"
$ws.Range("K7").Value = "2016-08-16 16:39:30"
$ws.Range("P7").Value = $errorDetail

$ws.Columns.Item(16).ColumnWidth = 39.17

# ---- de-de sheet ----
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("I7").Value = $displayName
$ws.Range("I7").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("I7"), $latestUrl, "", "", $displayName)

$ws.Range("J7").Value = "1f1910f2-2ea9-49f2-a4ad-68d0cddcb98c.3d92c1c223942b2c7e80af8189297e04f76f87e6.de-de.xlf"
$ws.Range("K7").Value = "2016-08-16 16:39:37"
$ws.Range("P7").Value = $errorDetail

$ws.Columns.Item(16).ColumnWidth = 39.17
